$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values like
# "0.999" or "622.81" are not silently coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '68.799.02'
$ws.Range('E2').Value = '  -1.19%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '3.628.32'
$ws.Range('E3').Value = '  -1.74%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '622.81'
$ws.Range('E5').Value = '  -7.40%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '157.63'
$ws.Range('E6').Value = '  -2.23%  '

$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.626.61'
$ws.Range('E7').Value = '  -1.73%  '

$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -1.99%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -2.27%  '

$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Value = '7.07'
$ws.Range('E11').Value = '  -0.73%  '

$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').Value = '  -1.81%  '

$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '0.0000226'
$ws.Range('E13').Value = '  -3.49%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.235.66'
$ws.Range('E14').Value = '  -1.88%  '

$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '31.96'
$ws.Range('E15').Value = '  -4.27%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.638.82'
$ws.Range('E16').Value = '  -0.91%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '68.884.53'
$ws.Range('E17').Value = '  -1.05%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.117'
$ws.Range('E18').Value = '  +0.66%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '6.44'
$ws.Range('E19').Value = '  -1.08%  '

$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '15.65'
$ws.Range('E20').Value = '  -3.34%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '10.15'
$ws.Range('E21').Value = '  +3.68%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '459.78'
$ws.Range('E22').Value = '  -2.77%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '0.642'
$ws.Range('E23').Value = '  -1.40%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '77.69'
$ws.Range('E24').Value = '  -3.11%  '

$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.782.09'
$ws.Range('E25').Value = '  -1.53%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0000118'
$ws.Range('E28').Value = '  -8.11%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '8.47'
$ws.Range('E29').Value = '  -7.42%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.59'
$ws.Range('E30').Value = '  -3.90%  '

$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.65'
$ws.Range('E31').Value = '  -4.17%  '

$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.11%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.94'
$ws.Range('E33').Value = '  -4.67%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '26.19'
$ws.Range('E34').Value = '  -2.67%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.162'
$ws.Range('E35').Value = '  -4.71%  '

$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.633.48'
$ws.Range('E36').Value = '  -1.47%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '6.28'
$ws.Range('E37').Value = '  -3.90%  '

$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '8.19'
$ws.Range('E38').Value = '  -3.94%  '

$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '178.59'
$ws.Range('E40').Value = '  +1.70%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  -0.32%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  -3.71%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '5.68'
$ws.Range('E43').Value = '  -7.19%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = '0.0886'
$ws.Range('E44').Value = '  -3.00%  '

$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.909'
$ws.Range('E45').Value = '  -3.08%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '29.03'
$ws.Range('E46').Value = '  +3.63%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '46.33'
$ws.Range('E47').Value = '  -1.43%  '

$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '2.64'
$ws.Range('E48').Value = '  -4.66%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '7.74'
$ws.Range('E49').Value = '  -2.02%  '

$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '0.000255'
$ws.Range('E50').Value = '  -7.60%  '

$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '1.02'
$ws.Range('E51').Value = '  -6.54%  '

# Restore the default (unstyled) cell format so the workbook styling
# matches the original — only cell contents should have changed.
$ws.Range("D2:E51").Style = "Normal"